# Automatische test-sync: 2025-08-06 19:44:50
#
# Adds the new "Inkoop / Bestellingen" mail-log entry (row 7 on "Logs"),
# rolls the new category count into the "Dashboard" summary table (row 3),
# widens the dashboard chart's category/value references to include the
# new row, and extends the conditional-formatting ranges on "Logs" from
# row 6 to row 7 so the new row picks up the same highlighting rules.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Logs" sheet - append the new mail entry as row 7
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A7").Value = "Bestel je 100 M5-bouten zodra je kan?"
$logs.Range("B7").Value = "mailmind.test@zohomail.eu"
$logs.Range("C7").Value = "Testmail #1: Bestel je 100 M5-bouten zodra je kan?"
$logs.Range("D7").Value = "Inkoop / Bestellingen"
$logs.Range("E7").Value = "Geachte klant,`nDank u voor uw e-mail. We willen u graag informeren dat wij geen bestelling hebben ontvangen voor 100 M5-bouten. Mocht dit een vergissing zijn, dan kunt u uw bestelling plaatsen via onze website of ons telefonisch bereiken.`nMet vriendelijke groet,`n[Naam Bedrijf]"
$logs.Range("F7").Value = "2025-08-06 19:44:13"
$logs.Range("G7").Value = "Ja"
$logs.Range("H7").Value = "Nee"
$logs.Range("I7").Value = "Ja"
$logs.Range("J7").Value = "Nee"

# The multi-line reply text in E7 makes the row auto-grow; re-fit it back
# down to the sheet's default height so row 7 stays un-customized, matching
# every other data row in this sheet.
$logs.Rows.Item(7).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 2. Extend the conditional-formatting ranges on "Logs" to cover row 7
# ---------------------------------------------------------------------------
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $newRange = $logs.Range($col + "2:" + $col + "7")
    $fc = $logs.Range($col + "2:" + $col + "6").FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------------
# 3. "Dashboard" sheet - add the new category/count row (row 3)
# ---------------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A3").Value = "Inkoop / Bestellingen"
$dashboard.Range("B3").Value = 1

# ---------------------------------------------------------------------------
# 4. Widen the dashboard chart's category/value series to include row 3
# ---------------------------------------------------------------------------
$chartObj = $dashboard.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection().Item(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$3"
$series.Values = "='Dashboard'!`$B`$2:`$B`$3"
